$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '41.146.19'
$ws.Range('E2').Value = '  -2.06%  '

# Row 3
$ws.Range('D3').Value = '2.177.75'
$ws.Range('E3').Value = '  -1.93%  '

# Row 4
$ws.Range('E4').Value = '  -0.18%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.18'
$ws.Range('E5').Value = '  -2.33%  '

# Row 6
$ws.Range('E6').Value = '  -0.99%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '70.40'
$ws.Range('E7').Value = '  -5.05%  '

# Row 8
$ws.Range('E8').Value = '  +0.14%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.581'
$ws.Range('E9').Value = '  -5.73%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.27'
$ws.Range('E10').Value = '  -8.71%  '

# Row 11
$ws.Range('E11').Value = '  -3.47%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.79'
$ws.Range('E12').Value = '  -5.29%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.101'
$ws.Range('E13').Value = '  -1.70%  '

# Row 14
$ws.Range('D14').Value = '2.502.28'
$ws.Range('E14').Value = '  -2.02%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.95'
$ws.Range('E15').Value = '  -2.46%  '

# Row 16
$ws.Range('E16').Value = '  -4.29%  '

# Row 17
$ws.Range('D17').Value = '2.184.68'
$ws.Range('E17').Value = '  -1.43%  '

# Row 18
$ws.Range('D18').Value = '41.071.34'
$ws.Range('E18').Value = '  -1.99%  '

# Row 19
$ws.Range('E19').Value = '  -7.36%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.50'
$ws.Range('E20').Value = '  -2.84%  '

# Row 21
$ws.Range('E21').Value = '  -4.01%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.04'
$ws.Range('E22').Value = '  -7.55%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '225.96'
$ws.Range('E23').Value = '  -1.97%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.96'
$ws.Range('E24').Value = '  -6.66%  '

# Row 25
$ws.Range('E25').Value = '  +0.12%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.91'
$ws.Range('E26').Value = '  -6.00%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.56'
$ws.Range('E27').Value = '  -0.94%  '

# Row 28
$ws.Range('E28').Value = '  -2.62%  '

# Row 29
$ws.Range('E29').Value = '  +1.83%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.31'
$ws.Range('E30').Value = '  +0.31%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.00'
$ws.Range('E31').Value = '  -3.13%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.13'
$ws.Range('E32').Value = '  +5.78%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0771'
$ws.Range('E33').Value = '  -3.99%  '

# Row 34
$ws.Range('E34').Value = '  -9.53%  '

# Row 35
$ws.Range('E35').Value = '  -3.11%  '

# Row 36
$ws.Range('E36').Value = '  -8.83%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.12'
$ws.Range('E37').Value = '  -4.35%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0286'
$ws.Range('E38').Value = '  -5.54%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.28'
$ws.Range('E39').Value = '  -5.98%  '

# Row 40
$ws.Range('E40').Value = '  -2.90%  '

# Row 41
$ws.Range('E41').Value = '  -4.06%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '60.34'
$ws.Range('E42').Value = '  -7.63%  '

# Row 43
$ws.Range('E43').Value = '  -4.68%  '

# Row 44
$ws.Range('E44').Value = '  -4.89%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0974'
$ws.Range('E45').Value = '  -3.74%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.50'
$ws.Range('E46').Value = '  -5.62%  '

# Row 47
$ws.Range('E47').Value = '  -2.83%  '

# Row 48
$ws.Range('E48').Value = '  -2.63%  '

# Row 49
$ws.Range('E49').Value = '  -8.42%  '

# Row 50
$ws.Range('E50').Value = '  -2.89%  '

# Row 51
$ws.Range('D51').Value = '2.380.13'
$ws.Range('E51').Value = '  -2.02%  '
